{"js": "// Office.js (Word JavaScript API) \u2014 body of: async (context) => { ... }\n//\n// Implements the \"post class 3 questions\" commit:\n//   1. Title: \"...preparing for January 21 class\" -> \"...preparing for Week 3 class\"\n//   2. \"What is an even study difference-in-difference\" -> \"...an event study...\"\n//   3. \"Appendix Table A2\" -> \"Appendix Table A3\"\n//   4. \"...last three paragraphs of section 6).\" -> \"...last three paragraphs on pg. 28).\"\n//   5. Remove the stray leftover \"_GoBack\" bookmark sitting mid-word in\n//      \"au|thors\" in the triple-difference question (no visible text change).\n\nasync function replaceText(context, searchText, newText, options) {\n  const opts = Object.assign({ matchCase: true, matchWholeWord: false }, options || {});\n  const results = context.document.body.search(searchText, opts);\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + searchText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. Title line: \"January 21\" -> \"Week 3\".\nawait replaceText(context, \"January 21\", \"Week 3\");\n\n// 2. \"an even study\" -> \"an event study\" (adds the missing \"t\").\nawait replaceText(context, \"an even study\", \"an event study\");\n\n// 3. \"Appendix Table A2\" -> \"Appendix Table A3\".\nawait replaceText(context, \"Appendix Table A2\", \"Appendix Table A3\");\n\n// 4. \"...last three paragraphs of section 6\" -> \"...last three paragraphs on pg. 28\".\nawait replaceText(\n  context,\n  \"last three paragraphs of section 6\",\n  \"last three paragraphs on pg. 28\"\n);\n\n// 5. Delete the orphaned \"_GoBack\" bookmark left behind mid-sentence, then\n//    rewrite the now-contiguous sentence so the two surrounding runs are\n//    merged back into a single run (matching the cleaned-up paragraph).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst finalSentence = \"third difference that the authors use in this study?\";\nawait replaceText(context, finalSentence, \"\\u0000TMP_MERGE_PLACEHOLDER\\u0000\");\nawait replaceText(context, \"\\u0000TMP_MERGE_PLACEHOLDER\\u0000\", finalSentence);\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word / $d (ActiveDocument) are pre-seeded by the harness.\n#\n# Implements the \"post class 3 questions\" commit:\n#   1. Title: \"...preparing for January 21 class\" -> \"...preparing for Week 3 class\"\n#   2. \"What is an even study difference-in-difference\" -> \"...an event study...\"\n#   3. \"Appendix Table A2\" -> \"Appendix Table A3\"\n#   4. \"...last three paragraphs of section 6).\" -> \"...last three paragraphs on pg. 28).\"\n#   5. Remove the stray leftover \"_GoBack\" bookmark sitting mid-word in\n#      \"au|thors\" in the triple-difference question (no visible text change).\n\n$d = $word.ActiveDocument\n\nfunction ReplaceOnce($searchText, $newText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.Text = $searchText\n    $find.MatchCase = $true\n    $found = $find.Execute()\n    if ($found) {\n        $rng.Text = $newText\n    } else {\n        throw \"Search text not found: $searchText\"\n    }\n}\n\n# 1. Title line: \"January 21\" -> \"Week 3\".\nReplaceOnce \"January 21\" \"Week 3\"\n\n# 2. \"an even study\" -> \"an event study\" (adds the missing \"t\").\nReplaceOnce \"an even study\" \"an event study\"\n\n# 3. \"Appendix Table A2\" -> \"Appendix Table A3\".\nReplaceOnce \"Appendix Table A2\" \"Appendix Table A3\"\n\n# 4. \"...last three paragraphs of section 6\" -> \"...last three paragraphs on pg. 28\".\nReplaceOnce \"last three paragraphs of section 6\" \"last three paragraphs on pg. 28\"\n\n# 5. Delete the orphaned \"_GoBack\" bookmark, then rewrite the now-contiguous\n#    sentence so the two surrounding runs merge back into a single run\n#    (matching the cleaned-up paragraph).\n$bm = $d.Bookmarks(\"_GoBack\")\n$bm.Delete()\n\n$finalSentence = \"third difference that the authors use in this study?\"\nReplaceOnce $finalSentence \"`0TMP_MERGE_PLACEHOLDER`0\"\nReplaceOnce \"`0TMP_MERGE_PLACEHOLDER`0\" $finalSentence\n"}
